# Update "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on row 4 of the
# zh-cn and de-de worksheets, as part of regenerating the handback
# report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-28 05:20:48"
$wsZhCn.Range("G4").Value = "2016-01-28 05:21:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-28 05:20:58"
$wsDeDe.Range("G4").Value = "2016-01-28 05:21:48"
